# Update supp table 2: reclassify specific MB, SHH samples as MB, SHH alpha
# based on genomic/expression evidence, on the "MB_SHH_subtype" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MB_SHH_subtype")

$rows = @(40, 45, 49, 50, 60, 63, 69, 81, 83, 103, 157, 171, 189, 214, 228, 244)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = "MB, SHH alpha"
    $ws.Range("I$r").Value = "Genomic/Expression"
}
